# Add 2022-Q4 data:
#  - a new worksheet "2022-Q4" (containing the single 006282 fund row) is
#    inserted right after the "总计" sheet, pushing every existing
#    2022-Q3 .. 2020-Q4 sheet one slot further along the tab strip.
#  - the "总计" (summary) sheet gets a new row 2 for the 2022-Q4 totals,
#    with every other row shifting down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right before the current
#    "2022-Q3" tab (i.e. right after "总计").
# ---------------------------------------------------------------------
$anchor = $wb.Worksheets.Item(2)
$q4 = $wb.Worksheets.Add($anchor)
$q4.Name = "2022-Q4"

# Header row (copy formatting from the "2022-Q3" sheet's header so the
# bold/centered/bordered style matches the other quarter sheets).
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Single data row: fund 006282.
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "006282"
$q4.Range("C2").Value = "上投摩根欧洲动力策略股票（QDII）"
$q4.Range("D2").Value = "0.43"
$q4.Range("E2").Value = "92.90"
$q4.Range("F2").Value = "2.67"
$q4.Range("G2").Value = "0.0115"
$q4.Range("H2").Value = 6

# Copy the header / first-column formatting from the "2022-Q3" sheet
# (now shifted to position 3) so the new tab matches the look of its
# siblings (bold header row + the A-column style used for row indices).
$src = $wb.Worksheets.Item(3)
$src.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$src.Range("A2").Copy()
$q4.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row 2 for 2022-Q4 and
#    shift the existing quarters down by one row.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$total.Rows.Item(2).Insert()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.01

# Restore formatting on the newly inserted row: column A keeps the
# centered/bordered index style, columns B:D stay plain (matching every
# other data row in the sheet).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B3:D3").Copy()
$total.Range("B2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. The "2020-Q4" tab (now the last sheet) stays the selected tab, same
#    as before the edit - activating the newly-inserted "2022-Q4" sheet
#    would otherwise steal the tabSelected flag.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
